$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric ("210.87", "2.50", etc.) must be
# pre-formatted as Text so Excel keeps them as literal strings instead of
# silently converting to a number (which would also eat formatting like
# trailing zeros, e.g. "2.50" -> 2.5).
$textCells = @("D5", "D6", "D10", "D11", "D14", "D15", "D16", "D20", "D21", "D22", "D24", "D25", "D27", "D29", "D30", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D49", "D51")
foreach ($r in $textCells) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.715.87"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.598.16"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "210.87"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.821.22"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.601.10"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "4.06"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "0.525"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "64.92"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "26.651.89"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "208.78"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "7.12"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "144.25"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "7.14"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +19.27%  "
$ws.Range("D35").Value = "1.278.94"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "2.50"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.598"
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.23"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.43"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "0.778"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "62.64"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "1.732.66"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "90.28"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0511"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0965"
$ws.Range("E50").Value = "  -9.24%  "
$ws.Range("D51").Value = "7.51"
